$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Cebollín baby" (Agrícola del
# Norte S.A. de Arica). It belongs right after the existing row 35, so
# insert a fresh row at position 36 — this pushes the former rows 36-96
# down to 37-97 (and the sheet's used range grows from R96 to R97), while
# the constant columns' formatting (e.g. the date style on column D)
# carries over automatically from the row below.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new record's data.
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 44792
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112038
$ws.Cells.Item(36, 7).Value = "Cebollín baby"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 250
$ws.Cells.Item(36, 11).Value = 3500
$ws.Cells.Item(36, 12).Value = 4000
$ws.Cells.Item(36, 13).Value = 3750
$ws.Cells.Item(36, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 1875
$ws.Cells.Item(36, 17).Value = 2
$ws.Cells.Item(36, 18).Value = "Hortaliza"
